# "Generate Report for Archive"
# - Status moves from "Ready for handoff" to "In Translation" (shared string
#   used by Overview!E2/F2 and the per-language sheets' Status column, C2).
# - The two language-status columns on Overview (E, F) and the Status column
#   on each language sheet (C) are narrowed to match the new report layout.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Update the status value everywhere it appears.
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"
$zhcn.Range("C2").Value = "In Translation"
$dede.Range("C2").Value = "In Translation"

# Narrow the columns that held the status text.
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5
$zhcn.Columns.Item(3).ColumnWidth = 12.5
$dede.Columns.Item(3).ColumnWidth = 12.5
